$d = $word.ActiveDocument

$replacements = @(
    @{old = "2025-11-08 Saturday"; new = "2025-11-09 Sunday"},
    @{old = "53÷5="; new = "55÷3="},
    @{old = "82÷5="; new = "66÷2="},
    @{old = "80÷7="; new = "75÷8="},
    @{old = "19÷5="; new = "14÷8="},
    @{old = "50÷5="; new = "35÷6="},
    @{old = "71÷7="; new = "55÷8="},
    @{old = "80÷9="; new = "77÷9="},
    @{old = "86÷9="; new = "61÷7="},
    @{old = "79÷4="; new = "99÷8="},
    @{old = "76÷7="; new = "32÷5="},
    @{old = "29÷6="; new = "77÷5="},
    @{old = "57÷3="; new = "97÷8="},
    @{old = "61÷8="; new = "47÷7="},
    @{old = "16÷6="; new = "93÷8="},
    @{old = "17÷2="; new = "64÷8="},
    @{old = "44÷9="; new = "85÷2="},
    @{old = "11÷4="; new = "75÷2="},
    @{old = "79÷6="; new = "77÷3="},
    @{old = "41÷3="; new = "43÷8="},
    @{old = "46÷9="; new = "18÷7="},
    @{old = "18÷2="; new = "31÷9="},
    @{old = "57÷7="; new = "11÷3="},
    @{old = "23÷7="; new = "54÷9="},
    @{old = "90÷2="; new = "71÷6="},
    @{old = "86÷7="; new = "35÷7="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
